$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at positions 461-462, shifting all existing
# data (old rows 461:565) down to 463:567.
$ws.Rows("461:462").Insert()

# --- New row 461 ---
$ws.Cells.Item(461,1).Value  = 11
$ws.Cells.Item(461,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(461,3).Value  = "Bíobío"
$ws.Cells.Item(461,4).Value  = 44798
$ws.Cells.Item(461,5).Value  = 8
$ws.Cells.Item(461,6).Value  = "Fruta"
$ws.Cells.Item(461,7).Value  = 100102
$ws.Cells.Item(461,8).Value  = "Cítricos"
$ws.Cells.Item(461,9).Value  = 100102003
$ws.Cells.Item(461,10).Value = "Limón"
$ws.Cells.Item(461,11).Value = "Sin especificar"
$ws.Cells.Item(461,12).Value = "1a amarillo"
$ws.Cells.Item(461,13).Value = 310
$ws.Cells.Item(461,14).Value = 5000
$ws.Cells.Item(461,15).Value = 5500
$ws.Cells.Item(461,16).Value = 5242
$ws.Cells.Item(461,17).Value = "$/malla 16 kilos"
$ws.Cells.Item(461,18).Value = "Región de O'Higgins"
$ws.Cells.Item(461,19).Value = 328
$ws.Cells.Item(461,20).Value = 16

# --- New row 462 ---
$ws.Cells.Item(462,1).Value  = 11
$ws.Cells.Item(462,2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(462,3).Value  = "Bíobío"
$ws.Cells.Item(462,4).Value  = 44798
$ws.Cells.Item(462,5).Value  = 8
$ws.Cells.Item(462,6).Value  = "Fruta"
$ws.Cells.Item(462,7).Value  = 100102
$ws.Cells.Item(462,8).Value  = "Cítricos"
$ws.Cells.Item(462,9).Value  = 100102003
$ws.Cells.Item(462,10).Value = "Limón"
$ws.Cells.Item(462,11).Value = "Sin especificar"
$ws.Cells.Item(462,12).Value = "2a amarillo"
$ws.Cells.Item(462,13).Value = 250
$ws.Cells.Item(462,14).Value = 4000
$ws.Cells.Item(462,15).Value = 4500
$ws.Cells.Item(462,16).Value = 4200
$ws.Cells.Item(462,17).Value = "$/malla 16 kilos"
$ws.Cells.Item(462,18).Value = "Región de O'Higgins"
$ws.Cells.Item(462,19).Value = 262
$ws.Cells.Item(462,20).Value = 16
